$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.733999999999999
$ws.Range("B7").Value = 5.344
$ws.Range("B16").Value = 5.339
$ws.Range("B28").Value = 6.034
$ws.Range("B29").Value = 5.615
$ws.Range("B32").Value = 6.845999999999999
$ws.Range("B40").Value = 9.327999999999999
$ws.Range("B52").Value = 5.359
$ws.Range("B57").Value = 5.090999999999999
$ws.Range("B66").Value = 5.114999999999999
$ws.Range("B100").Value = 5.931
